$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "27.426.99"
$ws.Range("D3").Value = "1.712.74"
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'224.57"
$ws.Range("E5").Value = "  -1.22%  "
$ws.Range("D6").Value = "'0.5333"
$ws.Range("E6").Value = "  -2.27%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'0.2671"
$ws.Range("E8").Value = "  -3.28%  "
$ws.Range("D9").Value = "'0.06616"
$ws.Range("E9").Value = "  -1.52%  "
$ws.Range("D10").Value = "'20.96"
$ws.Range("E10").Value = "  -4.55%  "
$ws.Range("D11").Value = "'0.07649"
$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.787.10"
$ws.Range("E12").Value = "  +2.78%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.573"
$ws.Range("E13").Value = "  -2.33%  "
$ws.Range("D14").Value = "1.948.53"
$ws.Range("E14").Value = "  -1.45%  "
$ws.Range("D15").Value = "'0.5782"
$ws.Range("E15").Value = "  -3.10%  "
$ws.Range("D16").Value = "0.0₅8197"
$ws.Range("E16").Value = "  -2.33%  "
$ws.Range("D17").Value = "'68.02"
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("D18").Value = "27.415.55"
$ws.Range("E18").Value = "  -1.30%  "
$ws.Range("D19").Value = "'217.29"
$ws.Range("E19").Value = "  -3.64%  "
$ws.Range("D20").Value = "'1.004"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "'4.679"
$ws.Range("E21").Value = "  -2.98%  "
$ws.Range("E22").Value = "  -3.70%  "
$ws.Range("D23").Value = "'5.988"
$ws.Range("E23").Value = "  -3.73%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "'142.63"
$ws.Range("E25").Value = "  -3.09%  "
$ws.Range("D26").Value = "'1.738"
$ws.Range("E26").Value = "  +2.12%  "
$ws.Range("E27").Value = "  -2.45%  "
$ws.Range("D28").Value = "'7.299"
$ws.Range("E28").Value = "  -2.04%  "
$ws.Range("D29").Value = "'16.33"
$ws.Range("E29").Value = "  -4.65%  "
$ws.Range("D30").Value = "'0.05419"
$ws.Range("E30").Value = "  -4.28%  "
$ws.Range("E31").Value = "  -1.44%  "
$ws.Range("D32").Value = "'3.517"
$ws.Range("E32").Value = "  -4.71%  "
$ws.Range("D33").Value = "'3.437"
$ws.Range("E33").Value = "  -2.06%  "
$ws.Range("D34").Value = "'1.651"
$ws.Range("E34").Value = "  -1.60%  "
$ws.Range("D35").Value = "'2.877"
$ws.Range("E35").Value = "  +0.81%  "
$ws.Range("D36").Value = "'0.9519"
$ws.Range("E36").Value = "  -2.53%  "
$ws.Range("D37").Value = "'2.416"
$ws.Range("E37").Value = "  -1.43%  "
$ws.Range("D38").Value = "'0.5880"
$ws.Range("E38").Value = "  -1.13%  "
$ws.Range("D39").Value = "'0.01636"
$ws.Range("E39").Value = "  -1.79%  "
$ws.Range("D40").Value = "'5.870"
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("D41").Value = "1.047.13"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").Value = "'1.005"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").Value = "'0.8415"
$ws.Range("E43").Value = "  -0.87%  "
$ws.Range("D44").Value = "'101.06"
$ws.Range("E44").Value = "  -0.74%  "
$ws.Range("D45").Value = "1.855.56"
$ws.Range("E45").Value = "  -1.45%  "
$ws.Range("E46").Value = "  -1.13%  "
$ws.Range("D47").Value = "'58.11"
$ws.Range("E47").Value = "  -1.80%  "
$ws.Range("E48").Value = "  +1.72%  "
$ws.Range("D49").Value = "'1.004"
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("D50").Value = "'8.093"
$ws.Range("E50").Value = "  -1.95%  "
$ws.Range("D51").Value = "'0.05240"
$ws.Range("E51").Value = "  -1.42%  "
